$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.486.85"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "3.691.19"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'681.59"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "4.315.78"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "'32.33"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "3.695.92"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "69.440.07"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").Value = "'16.00"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").Value = "'471.65"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "'0.651"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").Value = "'80.13"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").Value = "3.842.27"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").Value = "'10.90"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("D28").Value = "'9.13"
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'6.55"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'26.94"
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("D35").Value = "3.680.91"
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("E37").Value = "  +3.92%  "
$ws.Range("D38").Value = "'6.27"
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.28"
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("D43").Value = "'169.79"
$ws.Range("E43").Value = "  +2.88%  "
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").Value = "'47.08"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").Value = "'29.02"
$ws.Range("E46").Value = "  +4.69%  "
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("E48").Value = "  +3.60%  "
$ws.Range("E49").Value = "  +2.43%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  +0.17%  "
